# Auto-generated: update Leve profit-tracking cells (H:N) across sheets
# per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 6 (Leve Item ID 4564)
$ws_ALC.Cells.Item(6, 8).Value = 150
$ws_ALC.Cells.Item(6, 9).Value = 0
$ws_ALC.Cells.Item(6, 10).Value = 150
$ws_ALC.Cells.Item(6, 11).Value = 0
$ws_ALC.Cells.Item(6, 12).Value = 450
$ws_ALC.Cells.Item(6, 13).ClearContents()
$ws_ALC.Cells.Item(6, 14).Value = -674

# ALC row 9 (Leve Item ID 5487)
$ws_ALC.Cells.Item(9, 8).Value = 183.42857
$ws_ALC.Cells.Item(9, 9).Value = 71
$ws_ALC.Cells.Item(9, 10).Value = 333.33334
$ws_ALC.Cells.Item(9, 11).Value = 71
$ws_ALC.Cells.Item(9, 12).Value = 333.33334
$ws_ALC.Cells.Item(9, 13).Value = 98
$ws_ALC.Cells.Item(9, 14).Value = -671.33334

# ALC row 12 (Leve Item ID 5515)
$ws_ALC.Cells.Item(12, 8).Value = 1210
$ws_ALC.Cells.Item(12, 9).Value = 1251.8
$ws_ALC.Cells.Item(12, 10).Value = 1001
$ws_ALC.Cells.Item(12, 11).Value = 1251.8
$ws_ALC.Cells.Item(12, 12).Value = 1001
$ws_ALC.Cells.Item(12, 13).Value = -1081.8
$ws_ALC.Cells.Item(12, 14).Value = -1341

# ALC row 29 (Leve Item ID 4575)
$ws_ALC.Cells.Item(29, 8).Value = 3121.6667
$ws_ALC.Cells.Item(29, 9).Value = 1750
$ws_ALC.Cells.Item(29, 10).Value = 3396
$ws_ALC.Cells.Item(29, 11).Value = 5250
$ws_ALC.Cells.Item(29, 12).Value = 10188
$ws_ALC.Cells.Item(29, 13).Value = -4969
$ws_ALC.Cells.Item(29, 14).Value = -10750

# ALC row 38 (Leve Item ID 4599)
$ws_ALC.Cells.Item(38, 8).Value = 462.9091
$ws_ALC.Cells.Item(38, 9).Value = 227.42857
$ws_ALC.Cells.Item(38, 10).Value = 875
$ws_ALC.Cells.Item(38, 11).Value = 682.28571
$ws_ALC.Cells.Item(38, 12).Value = 2625
$ws_ALC.Cells.Item(38, 13).Value = -310.28571
$ws_ALC.Cells.Item(38, 14).Value = -3369

# ALC row 39 (Leve Item ID 4603)
$ws_ALC.Cells.Item(39, 8).Value = 100.333336
$ws_ALC.Cells.Item(39, 9).Value = 100.333336
$ws_ALC.Cells.Item(39, 11).Value = 301.000008
$ws_ALC.Cells.Item(39, 13).Value = -5.00000799999998

# ALC row 41 (Leve Item ID 5478)
$ws_ALC.Cells.Item(41, 8).Value = 391.66666
$ws_ALC.Cells.Item(41, 9).Value = 101.666664
$ws_ALC.Cells.Item(41, 10).Value = 488.33334
$ws_ALC.Cells.Item(41, 11).Value = 101.666664
$ws_ALC.Cells.Item(41, 12).Value = 488.33334
$ws_ALC.Cells.Item(41, 13).Value = 338.333336
$ws_ALC.Cells.Item(41, 14).Value = -1368.33334

# ALC row 58 (Leve Item ID 4606)
$ws_ALC.Cells.Item(58, 8).Value = 1734.2858
$ws_ALC.Cells.Item(58, 9).Value = 228
$ws_ALC.Cells.Item(58, 10).Value = 5500
$ws_ALC.Cells.Item(58, 11).Value = 684
$ws_ALC.Cells.Item(58, 12).Value = 16500
$ws_ALC.Cells.Item(58, 13).Value = -534
$ws_ALC.Cells.Item(58, 14).Value = -16800

# ALC row 74 (Leve Item ID 5507)
$ws_ALC.Cells.Item(74, 8).Value = 15629250
$ws_ALC.Cells.Item(74, 9).Value = 4003
$ws_ALC.Cells.Item(74, 10).Value = 17861428
$ws_ALC.Cells.Item(74, 11).Value = 4003
$ws_ALC.Cells.Item(74, 12).Value = 17861428
$ws_ALC.Cells.Item(74, 13).Value = -3067
$ws_ALC.Cells.Item(74, 14).Value = -17863300

# ALC row 77 (Leve Item ID 5507)
$ws_ALC.Cells.Item(77, 8).Value = 15629250
$ws_ALC.Cells.Item(77, 9).Value = 4003
$ws_ALC.Cells.Item(77, 10).Value = 17861428
$ws_ALC.Cells.Item(77, 11).Value = 20015
$ws_ALC.Cells.Item(77, 12).Value = 89307140
$ws_ALC.Cells.Item(77, 13).Value = -15335
$ws_ALC.Cells.Item(77, 14).Value = -89316500

# ALC row 87 (Leve Item ID 10651)
$ws_ALC.Cells.Item(87, 8).Value = 38604
$ws_ALC.Cells.Item(87, 10).Value = 38604
$ws_ALC.Cells.Item(87, 12).Value = 38604
$ws_ALC.Cells.Item(87, 14).Value = -41100

# ALC row 90 (Leve Item ID 10651)
$ws_ALC.Cells.Item(90, 8).Value = 38604
$ws_ALC.Cells.Item(90, 10).Value = 38604
$ws_ALC.Cells.Item(90, 12).Value = 115812
$ws_ALC.Cells.Item(90, 14).Value = -128292

# ALC row 92 (Leve Item ID 19901)
$ws_ALC.Cells.Item(92, 8).Value = 1574.375
$ws_ALC.Cells.Item(92, 9).Value = 1800
$ws_ALC.Cells.Item(92, 10).Value = 1198.3334
$ws_ALC.Cells.Item(92, 11).Value = 1800
$ws_ALC.Cells.Item(92, 12).Value = 1198.3334
$ws_ALC.Cells.Item(92, 13).Value = -552
$ws_ALC.Cells.Item(92, 14).Value = -3694.3334

# ALC row 97 (Leve Item ID 19885)
$ws_ALC.Cells.Item(97, 8).Value = 1102.5
$ws_ALC.Cells.Item(97, 10).Value = 1102.5
$ws_ALC.Cells.Item(97, 12).Value = 3307.5
$ws_ALC.Cells.Item(97, 14).Value = -4299.5

# ALC row 103 (Leve Item ID 19909)
$ws_ALC.Cells.Item(103, 8).Value = 250309.5
$ws_ALC.Cells.Item(103, 9).Value = 312755
$ws_ALC.Cells.Item(103, 11).Value = 938265
$ws_ALC.Cells.Item(103, 13).Value = -937679

# ALC row 106 (Leve Item ID 19903)
$ws_ALC.Cells.Item(106, 8).Value = 2605.5715
$ws_ALC.Cells.Item(106, 9).Value = 2085.8
$ws_ALC.Cells.Item(106, 11).Value = 2085.8
$ws_ALC.Cells.Item(106, 13).Value = -1454.8

# ALC row 107 (Leve Item ID 27766)
$ws_ALC.Cells.Item(107, 8).Value = 1490.625
$ws_ALC.Cells.Item(107, 9).Value = 1713.6364
$ws_ALC.Cells.Item(107, 11).Value = 1713.6364
$ws_ALC.Cells.Item(107, 13).Value = 206.3635999999999

# ALC row 113 (Leve Item ID 27775)
$ws_ALC.Cells.Item(113, 8).Value = 25003274
$ws_ALC.Cells.Item(113, 9).Value = 35716428
$ws_ALC.Cells.Item(113, 10).Value = 5916.5835
$ws_ALC.Cells.Item(113, 11).Value = 35716428
$ws_ALC.Cells.Item(113, 12).Value = 5916.5835
$ws_ALC.Cells.Item(113, 13).Value = -35713174
$ws_ALC.Cells.Item(113, 14).Value = -12424.5835

# ALC row 129 (Leve Item ID 36115)
$ws_ALC.Cells.Item(129, 8).Value = 851.1613
$ws_ALC.Cells.Item(129, 10).Value = 852.86664
$ws_ALC.Cells.Item(129, 12).Value = 2558.59992
$ws_ALC.Cells.Item(129, 14).Value = -12558.59992

# ALC row 132 (Leve Item ID 44049)
$ws_ALC.Cells.Item(132, 8).Value = 4150
$ws_ALC.Cells.Item(132, 9).Value = 4245.294
$ws_ALC.Cells.Item(132, 10).Value = 3745
$ws_ALC.Cells.Item(132, 11).Value = 12735.882
$ws_ALC.Cells.Item(132, 12).Value = 11235
$ws_ALC.Cells.Item(132, 13).Value = -10205.882
$ws_ALC.Cells.Item(132, 14).Value = -16295

# ALC row 137 (Leve Item ID 44013)
$ws_ALC.Cells.Item(137, 8).Value = 69594.13
$ws_ALC.Cells.Item(137, 9).Value = 5124.75
$ws_ALC.Cells.Item(137, 10).Value = 93037.55
$ws_ALC.Cells.Item(137, 11).Value = 15374.25
$ws_ALC.Cells.Item(137, 12).Value = 279112.65
$ws_ALC.Cells.Item(137, 13).Value = -12824.25
$ws_ALC.Cells.Item(137, 14).Value = -284212.65

# ALC row 138 (Leve Item ID 44169)
$ws_ALC.Cells.Item(138, 8).Value = 1860.5294
$ws_ALC.Cells.Item(138, 9).Value = 556.6
$ws_ALC.Cells.Item(138, 10).Value = 3114.3076
$ws_ALC.Cells.Item(138, 11).Value = 1669.8
$ws_ALC.Cells.Item(138, 12).Value = 9342.9228
$ws_ALC.Cells.Item(138, 13).Value = 3470.2
$ws_ALC.Cells.Item(138, 14).Value = -19622.9228

# ARM row 74 (Leve Item ID 44000)
$ws_ARM.Cells.Item(74, 8).Value = 2278.5
$ws_ARM.Cells.Item(74, 9).Value = 2262.9092
$ws_ARM.Cells.Item(74, 11).Value = 2262.9092
$ws_ARM.Cells.Item(74, 13).Value = -1388.9092

# ARM row 77 (Leve Item ID 44000)
$ws_ARM.Cells.Item(77, 8).Value = 2278.5
$ws_ARM.Cells.Item(77, 9).Value = 2262.9092
$ws_ARM.Cells.Item(77, 11).Value = 11314.546
$ws_ARM.Cells.Item(77, 13).Value = -6946.546

# ARM row 97 (Leve Item ID 19941)
$ws_ARM.Cells.Item(97, 8).Value = 2014.9412
$ws_ARM.Cells.Item(97, 9).Value = 1872.9166
$ws_ARM.Cells.Item(97, 10).Value = 2355.8
$ws_ARM.Cells.Item(97, 11).Value = 1872.9166
$ws_ARM.Cells.Item(97, 12).Value = 2355.8
$ws_ARM.Cells.Item(97, 13).Value = -1376.9166
$ws_ARM.Cells.Item(97, 14).Value = -3347.8

# ARM row 122 (Leve Item ID 36168)
$ws_ARM.Cells.Item(122, 8).Value = 1396.1875
$ws_ARM.Cells.Item(122, 9).Value = 1581.5834
$ws_ARM.Cells.Item(122, 10).Value = 840
$ws_ARM.Cells.Item(122, 11).Value = 4744.7502
$ws_ARM.Cells.Item(122, 12).Value = 2520
$ws_ARM.Cells.Item(122, 13).Value = -2294.7502
$ws_ARM.Cells.Item(122, 14).Value = -7420

# ARM row 132 (Leve Item ID 43997)
$ws_ARM.Cells.Item(132, 8).Value = 24611.912
$ws_ARM.Cells.Item(132, 9).Value = 2441.2856
$ws_ARM.Cells.Item(132, 11).Value = 7323.8568
$ws_ARM.Cells.Item(132, 13).Value = -4793.8568

# BSM row 94 (Leve Item ID 19939)
$ws_BSM.Cells.Item(94, 8).Value = 2756.2646
$ws_BSM.Cells.Item(94, 9).Value = 1437.3043
$ws_BSM.Cells.Item(94, 11).Value = 1437.3043
$ws_BSM.Cells.Item(94, 13).Value = -986.3043

# BSM row 134 (Leve Item ID 43998)
$ws_BSM.Cells.Item(134, 8).Value = 49445.59
$ws_BSM.Cells.Item(134, 9).Value = 67300.31
$ws_BSM.Cells.Item(134, 10).Value = 1833
$ws_BSM.Cells.Item(134, 11).Value = 201900.93
$ws_BSM.Cells.Item(134, 12).Value = 5499
$ws_BSM.Cells.Item(134, 13).Value = -199365.93
$ws_BSM.Cells.Item(134, 14).Value = -10569

# CRP row 58 (Leve Item ID 44021)
$ws_CRP.Cells.Item(58, 8).Value = 23124.87
$ws_CRP.Cells.Item(58, 9).Value = 1955.8889
$ws_CRP.Cells.Item(58, 10).Value = 36733.5
$ws_CRP.Cells.Item(58, 11).Value = 1955.8889
$ws_CRP.Cells.Item(58, 12).Value = 36733.5
$ws_CRP.Cells.Item(58, 13).Value = -1752.8889
$ws_CRP.Cells.Item(58, 14).Value = -37139.5

# CRP row 107 (Leve Item ID 27689)
$ws_CRP.Cells.Item(107, 8).Value = 1175.8572
$ws_CRP.Cells.Item(107, 9).Value = 1310
$ws_CRP.Cells.Item(107, 10).Value = 997
$ws_CRP.Cells.Item(107, 11).Value = 1310
$ws_CRP.Cells.Item(107, 12).Value = 997
$ws_CRP.Cells.Item(107, 13).Value = 610
$ws_CRP.Cells.Item(107, 14).Value = -4837

# CRP row 122 (Leve Item ID 36196)
$ws_CRP.Cells.Item(122, 8).Value = 836.55554
$ws_CRP.Cells.Item(122, 9).Value = 885.5714
$ws_CRP.Cells.Item(122, 10).Value = 665
$ws_CRP.Cells.Item(122, 11).Value = 2656.7142
$ws_CRP.Cells.Item(122, 12).Value = 1995
$ws_CRP.Cells.Item(122, 13).Value = -206.7142000000003
$ws_CRP.Cells.Item(122, 14).Value = -6895

# CRP row 134 (Leve Item ID 44020)
$ws_CRP.Cells.Item(134, 8).Value = 1196.1818
$ws_CRP.Cells.Item(134, 9).Value = 915.6923
$ws_CRP.Cells.Item(134, 10).Value = 1601.3334
$ws_CRP.Cells.Item(134, 11).Value = 2747.0769
$ws_CRP.Cells.Item(134, 12).Value = 4804.0002
$ws_CRP.Cells.Item(134, 13).Value = -212.0769
$ws_CRP.Cells.Item(134, 14).Value = -9874.0002

# CRP row 136 (Leve Item ID 44021)
$ws_CRP.Cells.Item(136, 8).Value = 23124.87
$ws_CRP.Cells.Item(136, 9).Value = 1955.8889
$ws_CRP.Cells.Item(136, 10).Value = 36733.5
$ws_CRP.Cells.Item(136, 11).Value = 5867.6667
$ws_CRP.Cells.Item(136, 12).Value = 110200.5
$ws_CRP.Cells.Item(136, 13).Value = -3317.6667
$ws_CRP.Cells.Item(136, 14).Value = -115300.5

# CUL row 126 (Leve Item ID 36045)
$ws_CUL.Cells.Item(126, 8).Value = 5436.6665
$ws_CUL.Cells.Item(126, 10).Value = 5436.6665
$ws_CUL.Cells.Item(126, 12).Value = 16309.9995
$ws_CUL.Cells.Item(126, 14).Value = -26189.9995

# CUL row 131 (Leve Item ID 36060)
$ws_CUL.Cells.Item(131, 8).Value = 764.21
$ws_CUL.Cells.Item(131, 10).Value = 791.0851
$ws_CUL.Cells.Item(131, 12).Value = 2373.2553
$ws_CUL.Cells.Item(131, 14).Value = -12453.2553

# GSM row 57 (Leve Item ID 2876)
$ws_GSM.Cells.Item(57, 8).Value = 29970
$ws_GSM.Cells.Item(57, 10).Value = 29970
$ws_GSM.Cells.Item(57, 12).Value = 29970
$ws_GSM.Cells.Item(57, 14).Value = -31610

# GSM row 97 (Leve Item ID 19940)
$ws_GSM.Cells.Item(97, 8).Value = 3137.7144
$ws_GSM.Cells.Item(97, 9).Value = 1558.75
$ws_GSM.Cells.Item(97, 10).Value = 8190.4
$ws_GSM.Cells.Item(97, 11).Value = 1558.75
$ws_GSM.Cells.Item(97, 12).Value = 8190.4
$ws_GSM.Cells.Item(97, 13).Value = -1062.75
$ws_GSM.Cells.Item(97, 14).Value = -9182.4

# GSM row 102 (Leve Item ID 36169)
$ws_GSM.Cells.Item(102, 8).Value = 1738.125
$ws_GSM.Cells.Item(102, 9).Value = 1782.619
$ws_GSM.Cells.Item(102, 11).Value = 1782.619
$ws_GSM.Cells.Item(102, 13).Value = -160.6189999999999

# GSM row 122 (Leve Item ID 36182)
$ws_GSM.Cells.Item(122, 8).Value = 2008.8334
$ws_GSM.Cells.Item(122, 9).Value = 2085.7144
$ws_GSM.Cells.Item(122, 10).Value = 1901.2
$ws_GSM.Cells.Item(122, 11).Value = 6257.1432
$ws_GSM.Cells.Item(122, 12).Value = 5703.6
$ws_GSM.Cells.Item(122, 13).Value = -3807.1432
$ws_GSM.Cells.Item(122, 14).Value = -10603.6

# LTW row 132 (Leve Item ID 44058)
$ws_LTW.Cells.Item(132, 8).Value = 2458.1428
$ws_LTW.Cells.Item(132, 9).Value = 1751.1
$ws_LTW.Cells.Item(132, 11).Value = 5253.299999999999
$ws_LTW.Cells.Item(132, 13).Value = -2723.299999999999

# LTW row 136 (Leve Item ID 44060)
$ws_LTW.Cells.Item(136, 8).Value = 126076.25
$ws_LTW.Cells.Item(136, 9).Value = 126076.25
$ws_LTW.Cells.Item(136, 11).Value = 378228.75
$ws_LTW.Cells.Item(136, 13).Value = -375678.75

# WVR row 122 (Leve Item ID 36208)
$ws_WVR.Cells.Item(122, 8).Value = 1829.85
$ws_WVR.Cells.Item(122, 9).Value = 1728.375
$ws_WVR.Cells.Item(122, 10).Value = 2235.75
$ws_WVR.Cells.Item(122, 11).Value = 5185.125
$ws_WVR.Cells.Item(122, 12).Value = 6707.25
$ws_WVR.Cells.Item(122, 13).Value = -2735.125
$ws_WVR.Cells.Item(122, 14).Value = -11607.25
